# Commit: "Change names from *img to img*"
#
# Renames the seven "*img" sheets to "img*" (himg->imgh, timg->imgt,
# simg->imgs, gimg->imgg, wimg->imgw, bimg->imgb, eimg->imge), and makes
# the "imge" sheet (formerly "eimg", the 17th/last sheet) the active tab
# instead of the "xbday" sheet (4th sheet) that was previously active.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Move the active/selected tab from "xbday" (index 3) to "imge" (index 16,
# the renamed former "eimg" sheet) - matches activeTab 3 -> 16 in the diff.
$wb.Worksheets.Item("imge").Activate()

# Minor incidental formatting touch-up on the "xbday" sheet that came along
# with the same save (E19:F19 pick up the same font formatting used
# elsewhere on the sheet).
$ws = $wb.Worksheets.Item("xbday")
$ws.Range("E19:F19").Font.Name = "Calibri"
